$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.145336747169495
$ws.Range("B1").Value = 2.126478672027588
$ws.Range("C1").Value = 10.18478012084961
$ws.Range("D1").Value = 2.534620046615601
$ws.Range("E1").Value = 1.280562281608582
